$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Part 3 description: append " and uses a complementary filter." right
#    after the existing "[Pitch, Roll]" run, as a brand-new run (so it does
#    not get merged into the preceding run's formatting/rsid).
# ---------------------------------------------------------------------------
$pitchRange = $d.Content
$found = $pitchRange.Find.Execute("[Pitch, Roll]", $false, $false, $false, $false, $false, $true, 1, $false)
if ($found) {
    $pitchRange.Collapse(0)
    $pitchRange.InsertAfter(" and uses a complementary filter.")
    $pitchRange.Font.NameAscii = "Times New Roman"
    $pitchRange.Font.Name = "Times New Roman"
    $pitchRange.Font.NameBi = "Times New Roman"
}

# ---------------------------------------------------------------------------
# 2) GitHub repo hyperlink: the visible URL text was split across three runs
#    ("https://github.com/nowei", "/", "cse562/tree/master/hw2"). Collapse
#    them back into a single run while keeping the Hyperlink character
#    style. There are two "https://github.com/nowei..." links in the doc
#    (the video link, then the repo link), so we anchor on the unique
#    "ub repo" text that immediately precedes the repo link's "GitHub
#    repo:" label, and search from there -- robust to offset drift caused
#    by the edit above. We then rewrite the tail text twice -- first to a
#    scratch value and then to the final value -- so the underlying
#    run-merge logic (which only triggers on an actual text change) fires
#    and really coalesces the runs, rather than silently no-op'ing on an
#    already-matching replace.
# ---------------------------------------------------------------------------
$docEnd = $d.Content.End
$labelRange = $d.Range(0, $docEnd)
$labelFound = $labelRange.Find.Execute("ub repo", $false, $false, $false, $false, $false, $true, 1, $false)
$anchor = $d.Range($labelRange.End, $docEnd)
$found = $anchor.Find.Execute("https://github.com/nowei", $false, $false, $false, $false, $false, $true, 1, $false)
if ($labelFound -and $found) {
    $tailStart = $anchor.End
    $tailEnd = $d.Content.End
    $tailProbe = $d.Range($tailStart, $tailEnd)
    $tailLen = ("/cse562/tree/master/hw2").Length
    $tailRange = $d.Range($tailStart, $tailStart + $tailLen)
    if ($tailRange.Text -eq "/cse562/tree/master/hw2") {
        $tailRange.Text = "/cse562/tree/master/hw2" + [char]8203
        $tailRange2 = $d.Range($tailStart, $tailRange.End)
        $tailRange2.Text = "/cse562/tree/master/hw2"
    }
}

Write-Output "done"
